$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "campus" column header after the existing "student_id" header
$ws.Range("H1").Value = "campus"

# Update selection to reflect the newly active cell (matches diff's sheetView selection)
$ws.Range("H1").Select()
